# Commit: "add files 12 sept"
# - Rename sheet "Sheet" -> "Sheet1"
# - Extend header row (1) from A1:D1 to A1:AA1 with new metric columns,
#   and apply bold + thin-box border + center/top aligned style to A1:AA1
# - Add a new data row (2) with results for the Qwen2.5-72B-Instruct run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rename the sheet -------------------------------------------------
$ws.Name = "Sheet1"

# --- header row (row 1) -------------------------------------------------
$headers = @{
    "A1" = "Model"
    "B1" = "Examples"
    "C1" = "OtherMetrics..."
    "D1" = "Time (s)"
    "E1" = "Date"
    "F1" = "Model Name"
    "G1" = "Exact Precision (Micro Avg)"
    "H1" = "Exact Recall (Micro Avg)"
    "I1" = "Exact F1 Score (Micro Avg)"
    "J1" = "Exact Precision (Macro Avg)"
    "K1" = "Exact Recall (Macro Avg)"
    "L1" = "Exact F1 Score (Macro Avg)"
    "M1" = "Exact Precision (Weighted Avg)"
    "N1" = "Exact Recall (Weighted Avg)"
    "O1" = "Exact F1 Score (Weighted Avg)"
    "P1" = "Partial Precision"
    "Q1" = "Partial Recall"
    "R1" = "Partial F1 Score"
    "S1" = "Partial TP"
    "T1" = "Partial FP"
    "U1" = "Partial FN"
    "V1" = "Support"
    "W1" = "Accuracy"
    "X1" = "Result Link"
    "Y1" = "Stats Link"
    "Z1" = "No of GPU Used"
    "AA1" = "Power Consumption"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# --- new data row (row 2) -----------------------------------------------
# Text value that looks like a date: prefix with an apostrophe so it is
# stored as a literal string instead of being parsed into a date serial
# number (mirrors typing '09/12/2025 into Excel).
$ws.Range("E2").Value = "'09/12/2025"

$ws.Range("F2").Value = "Qwen2.5-72B-Instruct"

$numericValues = @{
    "G2" = 0.3884615384615385
    "H2" = 0.3400673400673401
    "I2" = 0.3626570915619389
    "J2" = 0.2088235106282142
    "K2" = 0.1634986891141089
    "L2" = 0.1819014614836574
    "M2" = 0.4513397072480202
    "N2" = 0.3400673400673401
    "O2" = 0.3850494968929187
    "P2" = 0.4708171206225681
    "Q2" = 0.4087837837837838
    "R2" = 0.4376130198915009
    "S2" = 121
    "T2" = 136
    "U2" = 175
    "V2" = 297
    "W2" = 0.9460754943079689
}

foreach ($addr in $numericValues.Keys) {
    $ws.Range($addr).Value = $numericValues[$addr]
}

$ws.Range("X2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-72B-Instruct_4_shot.txt"
$ws.Range("Y2").Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-72B-Instruct_4_shot.txt"
$ws.Range("Z2").Value = "4 MLGPU"
$ws.Range("AA2").Value = "0.191 kWh"

# --- header styling (bold, thin box border, centered/top aligned) -------
# Build the combined format on a single cell first, then fan it out with a
# formats-only paste so the whole A1:AA1 row ends up sharing one style
# instead of Excel minting a fresh style record per incremental property
# write across the range.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160
$a1.Copy()
$ws.Range("A1:AA1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "applied header/data edits"
